# Update metadata & rebasing input files
# - Rebase the price-basis label used throughout the MNEMONIC table from
#   "Millions: 2019-20 prices" to "Millions: 2020-21 prices" (D7:D34).
# - Leave the user's selection on the range that was just updated (D7:D34),
#   matching where they had been working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7:D34").Value = "Millions: 2020-21 prices "

$ws.Range("D7:D34").Select()
